# ----------------------------------------------------------------------------
# This script applies two related changes to the IEEE13 workbook:
#
# 1. On the "Transformer" sheet, two new data rows are inserted into the
#    "Positive-Sequence 2W Transformer" table (rows 13 and 14), pushing the
#    "End of Positive-Sequence 2W Transformer" marker and every table below it
#    down by two rows.
#
# 2. On the "Bus" sheet, the three-phase bus entries (and the two-phase ones)
#    are re-sorted alphabetically by bus name (A, B, C) instead of the
#    previous C, A, B ordering. Only the bus name (column A) and the angle
#    (column E) change per row; the other columns stay associated with their
#    row position since they are identical across the group.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1) Transformer sheet: insert the two new Positive-Sequence 2W rows
# ----------------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("Transformer")

# Insert two blank rows at row 13, shifting the existing row 13 ("End of
# Positive-Sequence 2W Transformer") and everything below it down to row 15+.
$wsT.Rows.Item(13).Insert()
$wsT.Rows.Item(13).Insert()

# Row 13: "sub" transformer
$wsT.Cells.Item(13,1).Value = "sub"
$wsT.Cells.Item(13,2).Value = 1
$wsT.Cells.Item(13,3).Value = "sourcebus"
$wsT.Cells.Item(13,4).Value = "'650"
$wsT.Cells.Item(13,5).Value = 0.00001
$wsT.Cells.Item(13,6).Value = 0.00008000000000000001
$wsT.Cells.Item(13,7).Value = 0
$wsT.Cells.Item(13,8).Value = 0
$wsT.Cells.Item(13,9).Value = 1
$wsT.Cells.Item(13,10).Value = 1
$wsT.Cells.Item(13,11).Value = -30

# Row 14: "xfm1" transformer
$wsT.Cells.Item(14,1).Value = "xfm1"
$wsT.Cells.Item(14,2).Value = 1
$wsT.Cells.Item(14,3).Value = "xf1"
$wsT.Cells.Item(14,4).Value = "'634"
$wsT.Cells.Item(14,5).Value = 0.011
$wsT.Cells.Item(14,6).Value = 0.02
$wsT.Cells.Item(14,7).Value = 0
$wsT.Cells.Item(14,8).Value = 0
$wsT.Cells.Item(14,9).Value = 1
$wsT.Cells.Item(14,10).Value = 1
$wsT.Cells.Item(14,11).Value = 0

# ----------------------------------------------------------------------------
# 2) Bus sheet: re-sort the phase buses alphabetically (A, B, C)
# ----------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bus")

$ws.Cells.Item(3,1).Value = "632_A"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(4,1).Value = "632_B"
$ws.Cells.Item(4,5).Value = -120
$ws.Cells.Item(5,1).Value = "632_C"
$ws.Cells.Item(5,5).Value = 120
$ws.Cells.Item(6,1).Value = "633_A"
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(7,1).Value = "633_B"
$ws.Cells.Item(7,5).Value = -120
$ws.Cells.Item(8,1).Value = "633_C"
$ws.Cells.Item(8,5).Value = 120
$ws.Cells.Item(9,1).Value = "634_A"
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(10,1).Value = "634_B"
$ws.Cells.Item(10,5).Value = -120
$ws.Cells.Item(11,1).Value = "634_C"
$ws.Cells.Item(11,5).Value = 120
$ws.Cells.Item(12,1).Value = "645_B"
$ws.Cells.Item(12,5).Value = -120
$ws.Cells.Item(13,1).Value = "645_C"
$ws.Cells.Item(13,5).Value = 120
$ws.Cells.Item(14,1).Value = "646_B"
$ws.Cells.Item(14,5).Value = -120
$ws.Cells.Item(15,1).Value = "646_C"
$ws.Cells.Item(15,5).Value = 120
$ws.Cells.Item(16,1).Value = "650_A"
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(17,1).Value = "650_B"
$ws.Cells.Item(17,5).Value = -120
$ws.Cells.Item(18,1).Value = "650_C"
$ws.Cells.Item(18,5).Value = 120
$ws.Cells.Item(20,1).Value = "670_A"
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(21,1).Value = "670_B"
$ws.Cells.Item(21,5).Value = -120
$ws.Cells.Item(22,1).Value = "670_C"
$ws.Cells.Item(22,5).Value = 120
$ws.Cells.Item(23,1).Value = "671_A"
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(24,1).Value = "671_B"
$ws.Cells.Item(24,5).Value = -120
$ws.Cells.Item(25,1).Value = "671_C"
$ws.Cells.Item(25,5).Value = 120
$ws.Cells.Item(26,1).Value = "675_A"
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(27,1).Value = "675_B"
$ws.Cells.Item(27,5).Value = -120
$ws.Cells.Item(28,1).Value = "675_C"
$ws.Cells.Item(28,5).Value = 120
$ws.Cells.Item(29,1).Value = "680_A"
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(30,1).Value = "680_B"
$ws.Cells.Item(30,5).Value = -120
$ws.Cells.Item(31,1).Value = "680_C"
$ws.Cells.Item(31,5).Value = 120
$ws.Cells.Item(32,1).Value = "684_A"
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(33,1).Value = "684_C"
$ws.Cells.Item(33,5).Value = 120
$ws.Cells.Item(34,1).Value = "692_A"
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,1).Value = "692_B"
$ws.Cells.Item(35,5).Value = -120
$ws.Cells.Item(36,1).Value = "692_C"
$ws.Cells.Item(36,5).Value = 120
$ws.Cells.Item(37,1).Value = "brkr_A"
$ws.Cells.Item(37,5).Value = 0
$ws.Cells.Item(38,1).Value = "brkr_B"
$ws.Cells.Item(38,5).Value = -120
$ws.Cells.Item(39,1).Value = "brkr_C"
$ws.Cells.Item(39,5).Value = 120
$ws.Cells.Item(42,1).Value = "mid_A"
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,1).Value = "mid_B"
$ws.Cells.Item(43,5).Value = -120
$ws.Cells.Item(44,1).Value = "mid_C"
$ws.Cells.Item(44,5).Value = 120
$ws.Cells.Item(45,1).Value = "rg60_A"
$ws.Cells.Item(45,5).Value = 0
$ws.Cells.Item(46,1).Value = "rg60_B"
$ws.Cells.Item(46,5).Value = -120
$ws.Cells.Item(47,1).Value = "rg60_C"
$ws.Cells.Item(47,5).Value = 120
$ws.Cells.Item(48,1).Value = "sourcebus_A"
$ws.Cells.Item(48,5).Value = 0
$ws.Cells.Item(49,1).Value = "sourcebus_B"
$ws.Cells.Item(49,5).Value = -120
$ws.Cells.Item(50,1).Value = "sourcebus_C"
$ws.Cells.Item(50,5).Value = 120
$ws.Cells.Item(52,1).Value = "xf1_A"
$ws.Cells.Item(52,5).Value = 0
$ws.Cells.Item(53,1).Value = "xf1_B"
$ws.Cells.Item(53,5).Value = -120
$ws.Cells.Item(54,1).Value = "xf1_C"
$ws.Cells.Item(54,5).Value = 120
